$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transactions")
$ws.Activate()

# The first data row (old id 147, dated 42095, Accrual 289.73) is removed entirely;
# everything below shifts up one row.
$ws.Rows("2").Delete()

# What is now row 3 (old id 146, dated 41699, Accrual 19.23) is also removed entirely;
# everything below shifts up one row again.
$ws.Rows("3").Delete()

# A brand-new Accrual transaction row is inserted at row 4 (it inherits the
# formatting of the Accrual row directly above it).
$ws.Rows("4").Insert()

$ws.Cells.Item(4, 1).Value2 = 6885
$ws.Cells.Item(4, 2).Value = "Head Office"
$ws.Cells.Item(4, 3).Value2 = 41680
$ws.Cells.Item(4, 4).Value = "Accrual"
$ws.Cells.Item(4, 5).Value2 = 14.79
$ws.Cells.Item(4, 6).Value2 = 0
$ws.Cells.Item(4, 7).Value2 = 14.79
$ws.Cells.Item(4, 8).Value2 = 0
$ws.Cells.Item(4, 9).Value2 = 0
$ws.Cells.Item(4, 10).Value2 = 0

# Refresh the transaction IDs in column A for every remaining data row so
# they match the values produced by the corrected run.
$ws.Cells.Item(2, 1).Value2 = 6891
$ws.Cells.Item(3, 1).Value2 = 6886
$ws.Cells.Item(5, 1).Value2 = 6883
$ws.Cells.Item(6, 1).Value2 = 6884
$ws.Cells.Item(7, 1).Value2 = 6882

# Update the saved selection to match the new cursor position.
$ws.Range("D6").Select()
